$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.817903995513916
$ws.Range("B1").Value = 5.012856960296631
$ws.Range("C1").Value = 1.37090802192688
$ws.Range("D1").Value = 0.414430558681488
$ws.Range("E1").Value = 0.3257473111152649
